# Sprint05.xlsx -- "Closed out sprint 5, started sprint 6"
#
# Overview sheet: record the last two days of actual hours for the
# "Term Paper" tasks, close out the finished task, append a status note,
# and leave the selection where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Append a progress update to the sprint notes cell (A11) and keep its
# row height pinned at the original 15pt (it's a wrapped merged cell that
# would otherwise auto-grow to fit the extra line).
$ws.Range("A11").Value = "HW2 ate a lot of the last iteration. Going to try to focus on getting a draft of the paper done.`nGot some work done on paper, plan to have initial draft done during next iteration"
$ws.Rows.Item(11).RowHeight = 15

# Day J: actual hours logged came in at 1 (not the originally planned 2)
# for the three "Term Paper" tasks tracked that day.
$ws.Range("J4").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("J6").Value = 1

# Day K: the "Document findings" task wrapped up, so no hours remain.
$ws.Range("K10").Value = 0

# Leave the selection where the user was last working.
$ws.Range("K19").Select() | Out-Null
